$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "7323 ms"
$ws.Range("E3").Value = "7758 ms"
$ws.Range("E4").Value = "8483 ms"
$ws.Range("E5").Value = "5657 ms"
$ws.Range("E6").Value = "8345 ms"
$ws.Range("E7").Value = "5336 ms"
